$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H), mirroring the formatting of the existing header row
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data value for the new column
$ws.Range("H2").Value = 0
